$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)
$ws.Activate()
$ws.Range("G3").Value = 87
Write-Host "done"
